$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates -----------------------------------------------------
# F4: "Changing the position of the target" -> add "and get the point system working"
$ws.Range("F4").Value = "Changing the position of the target and get the point system working"

# G4: "Write the presentation" -> add "and record the presentation"
$ws.Range("G4").Value = "Write the presentation and record the presentation"

# New cell F3 gets the same (updated) "Changing the position..." text
$ws.Range("F3").Value = "Changing the position of the target and get the point system working"

# New cell F7 gets the same "Gather data on how it is performing" text as F6
$ws.Range("F7").Value = "Gather data on how it is performing"

# --- Style updates (apply "Check Cell" cell style, same as column B) --
$ws.Range("D4").Style = "Check Cell"
$ws.Range("E4").Style = "Check Cell"
$ws.Range("D6").Style = "Check Cell"
$ws.Range("E6").Style = "Check Cell"

# --- Column width updates ---------------------------------------------
$ws.Columns("F").ColumnWidth = 57.333333333333336
$ws.Columns("G").ColumnWidth = 42.166666666666664

# --- Selection change ---------------------------------------------------
$ws.Range("F10").Select()
